# feat: add 2022-Q1 data
#
# - Inserts a new sheet "2022-Q1" (with the per-fund breakdown, same layout
#   as the existing quarterly sheets) right before the "总计" summary sheet.
# - Adds a "2022-Q1" row at the top of the "总计" sheet's data table.

function Set-TextCell($ws, $row, $col, $val) {
    # Force a numeric-looking (or otherwise ambiguous) string to be stored
    # as text, matching the source inlineStr cells, then drop back to the
    # "Normal" cell style so we don't leave a stray NumberFormat-only style
    # behind (the source cells carry no explicit style here).
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value2 = $val
    $c.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Build the new "2022-Q1" sheet by cloning "2021-Q4" (identical header
#    row / column layout / cell styling) and dropping it just before 总计.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")
$template.Copy($totalSheet)

$ws = $wb.Worksheets.Item("2021-Q4 (2)")
$ws.Name = "2022-Q1"

# The template only has 2 data rows (rows 2-3); extend column A's styling
# down to row 7 so every data row (2-7) keeps the same "index" cell style.
$ws.Range("A2").Copy()
$ws.Range("A4:A7").PasteSpecial(-4122)

# -- row index column (0-based) --
$ws.Range("A2").Value2 = 0
$ws.Range("A3").Value2 = 1
$ws.Range("A4").Value2 = 2
$ws.Range("A5").Value2 = 3
$ws.Range("A6").Value2 = 4
$ws.Range("A7").Value2 = 5

# -- row 2: 005392 --
Set-TextCell $ws 2 2 "005392"
Set-TextCell $ws 2 3 "长信价值蓝筹两年定期开放灵活配置混合A"
Set-TextCell $ws 2 4 "5.17"
Set-TextCell $ws 2 5 "93.17"
Set-TextCell $ws 2 6 "7.80"
Set-TextCell $ws 2 7 "0.4033"
$ws.Cells.Item(2, 8).Value2 = 5

# -- row 3: 009911 --
Set-TextCell $ws 3 2 "009911"
Set-TextCell $ws 3 3 "长信价值蓝筹两年定期开放灵活配置混合C"
Set-TextCell $ws 3 4 "4.99"
Set-TextCell $ws 3 5 "93.17"
Set-TextCell $ws 3 6 "7.80"
Set-TextCell $ws 3 7 "0.3892"
$ws.Cells.Item(3, 8).Value2 = 5

# -- row 4: 010253 --
Set-TextCell $ws 4 2 "010253"
Set-TextCell $ws 4 3 "兴银中证500指数增强A"
Set-TextCell $ws 4 4 "2.19"
Set-TextCell $ws 4 5 "82.47"
Set-TextCell $ws 4 6 "0.94"
Set-TextCell $ws 4 7 "0.0206"
$ws.Cells.Item(4, 8).Value2 = 7

# -- row 5: 011205 --
Set-TextCell $ws 5 2 "011205"
Set-TextCell $ws 5 3 "兴银中证500指数增强C"
Set-TextCell $ws 5 4 "1.78"
Set-TextCell $ws 5 5 "82.47"
Set-TextCell $ws 5 6 "0.94"
Set-TextCell $ws 5 7 "0.0167"
$ws.Cells.Item(5, 8).Value2 = 7

# -- row 6: 004988 --
Set-TextCell $ws 6 2 "004988"
Set-TextCell $ws 6 3 "人保双利优选混合A"
Set-TextCell $ws 6 4 "0.58"
Set-TextCell $ws 6 5 "25.37"
Set-TextCell $ws 6 6 "0.58"
Set-TextCell $ws 6 7 "0.0034"
$ws.Cells.Item(6, 8).Value2 = 6

# -- row 7: 004989 (its 持有市值 rounds to a bare 0, stored as a number) --
Set-TextCell $ws 7 2 "004989"
Set-TextCell $ws 7 3 "人保双利优选混合C"
Set-TextCell $ws 7 4 "0.00"
Set-TextCell $ws 7 5 "25.37"
Set-TextCell $ws 7 6 "0.58"
$ws.Cells.Item(7, 7).Value2 = 0
$ws.Cells.Item(7, 8).Value2 = 6

# ---------------------------------------------------------------------
# 2. Prepend a "2022-Q1" row to the "总计" summary sheet, shifting the
#    existing three quarters down by one row.
# ---------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")

$tot.Range("B5").Value2 = $tot.Range("B4").Value2
$tot.Range("C5").Value2 = $tot.Range("C4").Value2
$tot.Range("D5").Value2 = $tot.Range("D4").Value2

$tot.Range("B4").Value2 = $tot.Range("B3").Value2
$tot.Range("C4").Value2 = $tot.Range("C3").Value2
$tot.Range("D4").Value2 = $tot.Range("D3").Value2

$tot.Range("B3").Value2 = $tot.Range("B2").Value2
$tot.Range("C3").Value2 = $tot.Range("C2").Value2
$tot.Range("D3").Value2 = $tot.Range("D2").Value2

# column A is just the 0-based row index; extend its style down to row 5
$tot.Range("A4").Copy()
$tot.Range("A5").PasteSpecial(-4122)
$tot.Range("A2").Value2 = 0
$tot.Range("A3").Value2 = 1
$tot.Range("A4").Value2 = 2
$tot.Range("A5").Value2 = 3

$tot.Range("B2").Value2 = "2022-Q1"
$tot.Range("C2").Value2 = 6
$tot.Range("D2").Value2 = 0.83

Write-Output "done"
